$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks up front (targets will be rebuilt below,
# row identities are changing so stale relationship targets must not survive)
$ws.Hyperlinks.Delete()

# Update data rows 2-18 with the freshly scraped listings
# row 2: 詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発
$ws.Range("A2").Value = "2025-09-06 01:14:13"
$ws.Range("B2").Value = "詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5377709"
$ws.Range("G2").Value = 245
$ws.Range("H2").Value = "🔥Next.js ◆開発,Node.js ◇アプリ"

# row 3: <Next.js、バックエンド開発> ガントチャートアプリの改修製造
$ws.Range("A3").Value = "2025-09-06 01:14:13"
$ws.Range("B3").Value = "<Next.js、バックエンド開発> ガントチャートアプリの改修製造"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5379158"
$ws.Range("G3").Value = 225
$ws.Range("H3").Value = "🔥Next.js ◆開発 ◇アプリ"

# row 4: 日本株・米国株ランキングメール自動配信システムの作成依頼。Pythonなど。
$ws.Range("A4").Value = "2025-09-06 01:14:13"
$ws.Range("B4").Value = "日本株・米国株ランキングメール自動配信システムの作成依頼。Pythonなど。"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5387291"
$ws.Range("G4").Value = 205
$ws.Range("H4").Value = "🔥Python"

# row 5: 【注目】公式LINEで診断機能を実現するGPT連動開発依頼
$ws.Range("A5").Value = "2025-09-06 01:14:13"
$ws.Range("B5").Value = "【注目】公式LINEで診断機能を実現するGPT連動開発依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5387629"
$ws.Range("G5").Value = 183
$ws.Range("H5").Value = "🔥GPT ◆開発"

# row 6: ポイントサイトのバックエンド開発:会員登録・ポイント付与/管理機能の開発【フルリ
$ws.Range("A6").Value = "2025-09-06 01:14:13"
$ws.Range("B6").Value = "ポイントサイトのバックエンド開発:会員登録・ポイント付与/管理機能の開発【フルリモート】"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5387645"
$ws.Range("G6").Value = 135
$ws.Range("H6").Value = "◆開発 ◇サイト"

# row 7: 【急募】解体工事見積書自動作成Excelシステム開発者募集
$ws.Range("A7").Value = "2025-09-06 01:14:13"
$ws.Range("B7").Value = "【急募】解体工事見積書自動作成Excelシステム開発者募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5387640"
$ws.Range("G7").Value = 113
$ws.Range("H7").Value = "◆開発,システム開発"

# row 8: 【LINEミニアプリ開発】紹介クーポン機能付きのミニアプリ開発依頼
$ws.Range("A8").Value = "2025-09-06 01:14:13"
$ws.Range("B8").Value = "【LINEミニアプリ開発】紹介クーポン機能付きのミニアプリ開発依頼"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5388066"
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = "◆開発 ◇アプリ"

# row 9: Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)
$ws.Range("A9").Value = "2025-09-06 01:14:13"
$ws.Range("B9").Value = "Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5379176"
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = "◆開発 ◇アプリ"

# row 10: 【急募】スキースノーボードスクール予約サイトの料金修正依頼
$ws.Range("A10").Value = "2025-09-06 01:14:13"
$ws.Range("B10").Value = "【急募】スキースノーボードスクール予約サイトの料金修正依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5374405"
$ws.Range("G10").Value = 38
$ws.Range("H10").Value = "◇サイト"

# row 11: 超初級・SE育成の技術研修 サブ講師
$ws.Range("A11").Value = "2025-09-06 01:14:13"
$ws.Range("B11").Value = "超初級・SE育成の技術研修 サブ講師"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5385021"
$ws.Range("G11").Value = 25
$ws.Range("H11").ClearContents()

# row 12: 【3Dシミュレーション】現場を自由に配置できるデータ制作依頼
$ws.Range("A12").Value = "2025-09-06 01:14:13"
$ws.Range("B12").Value = "【3Dシミュレーション】現場を自由に配置できるデータ制作依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5387676"
$ws.Range("G12").Value = 25
$ws.Range("H12").ClearContents()

# row 13: 【SRE / インフラエンジニア募集(基本リモート・金融系案件)】
$ws.Range("A13").Value = "2025-09-06 01:14:13"
$ws.Range("B13").Value = "【SRE / インフラエンジニア募集(基本リモート・金融系案件)】"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5371075"
$ws.Range("G13").Value = 25
$ws.Range("H13").ClearContents()

# row 14: 限定公開 PR 限定公開の仕事
$ws.Range("A14").Value = "2025-09-06 01:14:13"
$ws.Range("B14").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5385681"
$ws.Range("G14").Value = 25
$ws.Range("H14").ClearContents()

# row 15: 限定公開 限定公開の仕事
$ws.Range("A15").Value = "2025-09-06 01:14:13"
$ws.Range("B15").Value = "限定公開 限定公開の仕事"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5387827"
$ws.Range("G15").Value = 18
$ws.Range("H15").ClearContents()

# row 16: 【NAS導入】VPN設定とネットワークドライブの構築支援
$ws.Range("A16").Value = "2025-09-06 01:14:13"
$ws.Range("B16").Value = "【NAS導入】VPN設定とネットワークドライブの構築支援"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5387417"
$ws.Range("G16").Value = 13
$ws.Range("H16").ClearContents()

# row 17: MT4 ea 制作
$ws.Range("A17").Value = "2025-09-06 01:14:13"
$ws.Range("B17").Value = "MT4 ea 制作"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5387933"
$ws.Range("G17").Value = 10
$ws.Range("H17").ClearContents()

# row 18: 【急募】Excelで自動シート増加と画像トリミングを実現!
$ws.Range("A18").Value = "2025-09-06 01:14:13"
$ws.Range("B18").Value = "【急募】Excelで自動シート増加と画像トリミングを実現!"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5387258"
$ws.Range("G18").Value = 10
$ws.Range("H18").ClearContents()

# Rows 19-32 no longer exist in the refreshed listing -> remove them entirely
$ws.Range("A19:A32").EntireRow.Delete()

# Re-create the hyperlinks for the URL column (F2:F18) with correct targets
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5377709", "", "", "https://www.lancers.jp/work/detail/5377709")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5379158", "", "", "https://www.lancers.jp/work/detail/5379158")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5387291", "", "", "https://www.lancers.jp/work/detail/5387291")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5387629", "", "", "https://www.lancers.jp/work/detail/5387629")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5387645", "", "", "https://www.lancers.jp/work/detail/5387645")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5387640", "", "", "https://www.lancers.jp/work/detail/5387640")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5388066", "", "", "https://www.lancers.jp/work/detail/5388066")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5379176", "", "", "https://www.lancers.jp/work/detail/5379176")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5374405", "", "", "https://www.lancers.jp/work/detail/5374405")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5385021", "", "", "https://www.lancers.jp/work/detail/5385021")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5387676", "", "", "https://www.lancers.jp/work/detail/5387676")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5371075", "", "", "https://www.lancers.jp/work/detail/5371075")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5385681", "", "", "https://www.lancers.jp/work/detail/5385681")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5387827", "", "", "https://www.lancers.jp/work/detail/5387827")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5387417", "", "", "https://www.lancers.jp/work/detail/5387417")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5387933", "", "", "https://www.lancers.jp/work/detail/5387933")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5387258", "", "", "https://www.lancers.jp/work/detail/5387258")
